$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.118.73"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.631.78"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.10"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.13"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "1.859.63"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "1.632.69"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.69"
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").Value = "27.090.95"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.07"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.39"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.18"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.58"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").Value = "1.301.54"
$ws.Range("E34").Value = "  +2.29%  "
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  +2.11%  "
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").Value = "1.768.83"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.66"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.820"
$ws.Range("E49").Value = "  +22.07%  "
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.57"
$ws.Range("E51").Value = "  -1.75%  "
